$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on columns/cells that must stay as text
# (so numeric-looking and date-looking strings aren't auto-converted)
$ws.Range("B1:E6").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "ProductName"
$ws.Range("C1").Value = "LOT"
$ws.Range("D1").Value = "DateReceived"
$ws.Range("E1").Value = "DateReceivedIni"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "DemoProd1"
$ws.Range("C2").Value = "LOT-001"
$ws.Range("D2").Value = "2026-02-02"
$ws.Range("E2").Value = "20260202"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "3"
$ws.Range("C3").Value = "4"
$ws.Range("D3").Value = "5"
$ws.Range("E3").Value = "6"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "3"
$ws.Range("C4").Value = "3"
$ws.Range("D4").Value = "3"
$ws.Range("E4").Value = "3"

# Row 5
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "6"
$ws.Range("C5").Value = "6"
$ws.Range("D5").Value = "6"
$ws.Range("E5").Value = "6"

# Row 6
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Hey"
$ws.Range("C6").Value = "1"
$ws.Range("D6").Value = "1"
$ws.Range("E6").Value = "1"

# Remove the temporary "Text" number format now that the text values are
# locked in, so cells don't carry a leftover style index.
$ws.Range("B1:E6").ClearFormats()
